$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table1")

# Add the new login row (row 2) to the LoginInfo table.
# Write Username/Email first, then Name, then the numeric Password so the
# shared-string table ends up with the same ordering the workbook author's
# Excel session produced: Yam, afasfasf, "Yam biton".
$ws.Cells.Item(2, 2).Value = "Yam"
$ws.Cells.Item(2, 3).Value = "afasfasf"
$ws.Cells.Item(2, 1).Value = "Yam biton"
$ws.Cells.Item(2, 4).Value = 123456789

# Move the active selection down to B4, matching the saved view state.
$ws.Range("B4").Select()

# Switch the workbook's body font from Calibri to Arial.
$normal = $wb.Styles.Item("Normal")
$normal.Font.Name = "Arial"
